{"js": "// Update the SVAROM\u00c5L filing from the \"muntlig f\u00f6rberedelse\" draft to the\n// final pre-\"huvudf\u00f6rhandling\" version: refresh the cover-page date line,\n// split the status table's last row into a completed \"Bevisning inl\u00e4mnad\"\n// row plus a new \"Huvudf\u00f6rhandling\" row, mark the Arabic-evidence\n// translations as done, note that all evidence was filed in time, and\n// refresh the signature date.\n\nconst body = context.document.body;\n\n// 1) Sub-title under \"M\u00e5l nr T 4438-25\".\nconst title = body.search(\n  \"Uppdaterad inlaga efter muntlig f\u00f6rberedelse 2026-02-16\",\n  { matchCase: true }\n);\n\n// 2) Cover \"Datum:\" line.\nconst datum = body.search(\"Datum: 2026-02-16\", { matchCase: true });\n\n// 3) OBS note about Arabic evidence translations.\nconst obs = body.search(\n  \"OBS: All arabisk bevisning kommer att f\u00f6rses med auktoriserad \u00f6vers\u00e4ttning enligt r\u00e4ttens instruktioner fr\u00e5n 2026-02-16.\",\n  { matchCase: true }\n);\n\n// 4) Reservation-of-right-to-invoke-further-evidence sentence.\nconst reserve = body.search(\n  \"Mohammad och Joumana f\u00f6rbeh\u00e5ller sig r\u00e4tten att \u00e5beropa ytterligare bevisning f\u00f6re fristen 2026-03-20.\",\n  { matchCase: true }\n);\n\n// 5) Signature line at the bottom of the document.\nconst sign = body.search(\"Eskilstuna den 2026-02-16\", { matchCase: true });\n\nawait context.sync();\n\ntitle.items[0].insertText(\n  \"Slutlig inlaga inf\u00f6r huvudf\u00f6rhandling 2026-05-18\",\n  \"Replace\"\n);\ndatum.items[0].insertText(\"Datum: 2026-05-18\", \"Replace\");\nobs.items[0].insertText(\n  \"Auktoriserade \u00f6vers\u00e4ttningar av all arabisk bevisning har bifogats i enlighet med r\u00e4ttens instruktioner fr\u00e5n 2026-02-16.\",\n  \"Replace\"\n);\nreserve.items[0].insertText(\n  \"All bevisning har inl\u00e4mnats inom fristen 2026-03-20.\",\n  \"Replace\"\n);\nsign.items[0].insertText(\"Eskilstuna den 2026-05-18\", \"Replace\");\n\n// 6) Status table (first table in the body): its last row used to read\n//    \"N\u00e4sta f\u00f6rhandling | 2026-05-18\". Re-purpose it as \"Bevisning\n//    inl\u00e4mnad | 2026-03-20\" and append a new \"Huvudf\u00f6rhandling |\n//    2026-05-18\" row after it.\nconst table = body.tables.getFirst();\ntable.getCell(4, 0).value = \"Bevisning inl\u00e4mnad\";\ntable.getCell(4, 1).value = \"2026-03-20\";\ntable.addRows(\"End\", 1, [[\"Huvudf\u00f6rhandling\", \"2026-05-18\"]]);\n\nawait context.sync();\n", "ps1": "# Update the SVAROM\u00c5L filing from the \"muntlig f\u00f6rberedelse\" draft to the\n# final pre-\"huvudf\u00f6rhandling\" version: refresh the cover-page date line,\n# split the status table's last row into a completed \"Bevisning inl\u00e4mnad\"\n# row plus a new \"Huvudf\u00f6rhandling\" row, mark the Arabic-evidence\n# translations as done, note that all evidence was filed in time, and\n# refresh the signature date.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# 1) Sub-title under \"M\u00e5l nr T 4438-25\".\nReplace-Text \"Uppdaterad inlaga efter muntlig f\u00f6rberedelse 2026-02-16\" \"Slutlig inlaga inf\u00f6r huvudf\u00f6rhandling 2026-05-18\"\n\n# 2) Cover \"Datum:\" line.\nReplace-Text \"Datum: 2026-02-16\" \"Datum: 2026-05-18\"\n\n# 3) OBS note about Arabic evidence translations.\nReplace-Text \"OBS: All arabisk bevisning kommer att f\u00f6rses med auktoriserad \u00f6vers\u00e4ttning enligt r\u00e4ttens instruktioner fr\u00e5n 2026-02-16.\" \"Auktoriserade \u00f6vers\u00e4ttningar av all arabisk bevisning har bifogats i enlighet med r\u00e4ttens instruktioner fr\u00e5n 2026-02-16.\"\n\n# 4) Reservation-of-right-to-invoke-further-evidence sentence.\nReplace-Text \"Mohammad och Joumana f\u00f6rbeh\u00e5ller sig r\u00e4tten att \u00e5beropa ytterligare bevisning f\u00f6re fristen 2026-03-20.\" \"All bevisning har inl\u00e4mnats inom fristen 2026-03-20.\"\n\n# 5) Signature line at the bottom of the document.\nReplace-Text \"Eskilstuna den 2026-02-16\" \"Eskilstuna den 2026-05-18\"\n\n# 6) Status table (first table in the body): its last row used to read\n#    \"N\u00e4sta f\u00f6rhandling | 2026-05-18\". Re-purpose it as \"Bevisning\n#    inl\u00e4mnad | 2026-03-20\" and append a new \"Huvudf\u00f6rhandling |\n#    2026-05-18\" row after it.\n$t = $d.Tables.Item(1)\n$lastRow = $t.Rows.Count\n$t.Cell($lastRow, 1).Range.Text = \"Bevisning inl\u00e4mnad\"\n$t.Cell($lastRow, 2).Range.Text = \"2026-03-20\"\n\n$newRow = $t.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"Huvudf\u00f6rhandling\"\n$newRow.Cells.Item(2).Range.Text = \"2026-05-18\"\n"}
